$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lũy kế tháng SÓC TRĂNG")

# Update the "last_edited_time" text value (shared string) used by several
# rows (2,3,6,7,8,11,13) from 2024-08-03T03:54:00.000Z to
# 2024-08-03T20:14:00.000Z
$ws.Range("D2").Value = "2024-08-03T20:14:00.000Z"
$ws.Range("D3").Value = "2024-08-03T20:14:00.000Z"
$ws.Range("D6").Value = "2024-08-03T20:14:00.000Z"
$ws.Range("D7").Value = "2024-08-03T20:14:00.000Z"
$ws.Range("D8").Value = "2024-08-03T20:14:00.000Z"
$ws.Range("D11").Value = "2024-08-03T20:14:00.000Z"
$ws.Range("D13").Value = "2024-08-03T20:14:00.000Z"

# Update the numeric amounts for row 7 (Tháng 8 record)
$ws.Range("S7").Value = 1560000     # properties.Chi tiêu.number
$ws.Range("W7").Value = 10440000    # properties.Lũy kế.formula.number
$ws.Range("AE7").Value = 12000000   # properties.Tổng doanh thu.formula.number
$ws.Range("AH7").Value = 12000000   # properties.Đã thanh toán.number
$ws.Range("AQ7").Value = 14000000   # properties.Đơn giá.number
